$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 358.125
$ws.Range("I6").Value = 378.66666
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 1135.99998
$ws.Range("L6").Value = 150
$ws.Range("M6").Value = -1023.99998
$ws.Range("N6").Value = -374
$ws.Range("H20").Value = 1210.5
$ws.Range("I20").Value = 1210.5
$ws.Range("K20").Value = 1210.5
$ws.Range("M20").Value = -980.5
$ws.Range("H28").Value = 806.5
$ws.Range("I28").Value = 788.9091
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 788.9091
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = -303.9091
$ws.Range("N28").Value = -1970
$ws.Range("H35").Value = 1210.5
$ws.Range("I35").Value = 1210.5
$ws.Range("K35").Value = 1210.5
$ws.Range("M35").Value = -831.5
$ws.Range("H39").Value = 197.92857
$ws.Range("I39").Value = 37
$ws.Range("K39").Value = 111
$ws.Range("M39").Value = 185
$ws.Range("H42").Value = 89
$ws.Range("I42").Value = 89
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 267
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -37
$ws.Range("N42").ClearContents()
$ws.Range("H76").Value = 7660.2256
$ws.Range("I76").Value = 10332.066
$ws.Range("J76").Value = 5155.375
$ws.Range("K76").Value = 10332.066
$ws.Range("L76").Value = 5155.375
$ws.Range("M76").Value = -10017.066
$ws.Range("N76").Value = -5785.375
$ws.Range("H79").Value = 7660.2256
$ws.Range("I79").Value = 10332.066
$ws.Range("J79").Value = 5155.375
$ws.Range("K79").Value = 10332.066
$ws.Range("L79").Value = 5155.375
$ws.Range("M79").Value = -9240.066000000001
$ws.Range("N79").Value = -7339.375
$ws.Range("H106").Value = 3180.8462
$ws.Range("I106").Value = 2570.7058
$ws.Range("J106").Value = 4333.3335
$ws.Range("K106").Value = 2570.7058
$ws.Range("L106").Value = 4333.3335
$ws.Range("M106").Value = -1939.7058
$ws.Range("N106").Value = -5595.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8520.705
$ws.Range("I32").Value = 4416.311
$ws.Range("J32").Value = 84452
$ws.Range("K32").Value = 4416.311
$ws.Range("L32").Value = 84452
$ws.Range("M32").Value = -4129.311
$ws.Range("N32").Value = -85026
$ws.Range("H34").Value = 4000
$ws.Range("J34").Value = 4000
$ws.Range("L34").Value = 4000
$ws.Range("N34").Value = -4542
$ws.Range("H61").Value = 2933.3438
$ws.Range("I61").Value = 2220.3044
$ws.Range("J61").Value = 4755.5557
$ws.Range("K61").Value = 2220.3044
$ws.Range("L61").Value = 4755.5557
$ws.Range("M61").Value = -2008.3044
$ws.Range("N61").Value = -5179.5557
$ws.Range("H74").Value = 1637.4073
$ws.Range("I74").Value = 969.093
$ws.Range("J74").Value = 4249.909
$ws.Range("K74").Value = 969.093
$ws.Range("L74").Value = 4249.909
$ws.Range("M74").Value = -95.09299999999996
$ws.Range("N74").Value = -5997.909
$ws.Range("H77").Value = 1637.4073
$ws.Range("I77").Value = 969.093
$ws.Range("J77").Value = 4249.909
$ws.Range("K77").Value = 4845.465
$ws.Range("L77").Value = 21249.545
$ws.Range("M77").Value = -477.4650000000001
$ws.Range("N77").Value = -29985.545
$ws.Range("H110").Value = 1168.65
$ws.Range("I110").Value = 1010.17645
$ws.Range("K110").Value = 1010.17645
$ws.Range("M110").Value = 1034.82355
$ws.Range("H132").Value = 3649.6
$ws.Range("I132").Value = 3695.68
$ws.Range("J132").Value = 3572.8
$ws.Range("K132").Value = 11087.04
$ws.Range("L132").Value = 10718.4
$ws.Range("M132").Value = -8557.039999999999
$ws.Range("N132").Value = -15778.4
$ws.Range("H136").Value = 2933.3438
$ws.Range("I136").Value = 2220.3044
$ws.Range("J136").Value = 4755.5557
$ws.Range("K136").Value = 6660.9132
$ws.Range("L136").Value = 14266.6671
$ws.Range("M136").Value = -4110.9132
$ws.Range("N136").Value = -19366.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 20000
$ws.Range("J62").Value = 20000
$ws.Range("L62").Value = 20000
$ws.Range("N62").Value = -21372
$ws.Range("H65").Value = 20000
$ws.Range("J65").Value = 20000
$ws.Range("L65").Value = 60000
$ws.Range("N65").Value = -66864
$ws.Range("H105").Value = 2272.2708
$ws.Range("I105").Value = 2133.3948
$ws.Range("J105").Value = 2800
$ws.Range("K105").Value = 2133.3948
$ws.Range("L105").Value = 2800
$ws.Range("M105").Value = -386.3948
$ws.Range("N105").Value = -6294
$ws.Range("H134").Value = 28708.854
$ws.Range("I134").Value = 32956.453
$ws.Range("K134").Value = 98869.359
$ws.Range("M134").Value = -96334.359

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H31").Value = 3349.3333
$ws.Range("I31").Value = 2168.56
$ws.Range("J31").Value = 5457.857
$ws.Range("K31").Value = 2168.56
$ws.Range("L31").Value = 5457.857
$ws.Range("N31").Value = -6047.857
$ws.Range("M31").Value = -1873.56
$ws.Range("H34").Value = 3349.3333
$ws.Range("I34").Value = 2168.56
$ws.Range("J34").Value = 5457.857
$ws.Range("K34").Value = 2168.56
$ws.Range("L34").Value = 5457.857
$ws.Range("N34").Value = -5861.857
$ws.Range("M34").Value = -1966.56
$ws.Range("H47").Value = 16535.5
$ws.Range("J47").Value = 16535.5
$ws.Range("L47").Value = 16535.5
$ws.Range("N47").Value = -17667.5
$ws.Range("H51").Value = 22240
$ws.Range("J51").Value = 22240
$ws.Range("L51").Value = 22240
$ws.Range("N51").Value = -23712
$ws.Range("H59").Value = 55706.715
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 55706.715
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 55706.715
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -57996.715
$ws.Range("H61").Value = 22240
$ws.Range("J61").Value = 22240
$ws.Range("L61").Value = 22240
$ws.Range("N61").Value = -22936
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H132").Value = 2132.204
$ws.Range("I132").Value = 1134.9166
$ws.Range("J132").Value = 3089.6
$ws.Range("K132").Value = 3404.7498
$ws.Range("L132").Value = 9268.799999999999
$ws.Range("M132").Value = -874.7498000000001
$ws.Range("N132").Value = -14328.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 2275
$ws.Range("J42").Value = 2275
$ws.Range("L42").Value = 6825
$ws.Range("N42").Value = -7893
$ws.Range("H48").Value = 3872.6
$ws.Range("I48").Value = 799
$ws.Range("J48").Value = 4641
$ws.Range("K48").Value = 2397
$ws.Range("L48").Value = 13923
$ws.Range("M48").Value = -2147
$ws.Range("N48").Value = -14423
$ws.Range("H55").Value = 3362.5
$ws.Range("J55").Value = 3728.5715
$ws.Range("L55").Value = 11185.7145
$ws.Range("N55").Value = -11539.7145
$ws.Range("H56").Value = 3978.5715
$ws.Range("I56").Value = 3978.5715
$ws.Range("K56").Value = 3978.5715
$ws.Range("M56").Value = -3448.5715
$ws.Range("H131").Value = 771.2963
$ws.Range("I131").Value = 376.53333
$ws.Range("J131").Value = 923.12823
$ws.Range("K131").Value = 1129.59999
$ws.Range("L131").Value = 2769.38469
$ws.Range("M131").Value = 3910.40001
$ws.Range("N131").Value = -12849.38469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 8000
$ws.Range("J48").Value = 8000
$ws.Range("L48").Value = 8000
$ws.Range("N48").Value = -8970
$ws.Range("H132").Value = 2685.2144
$ws.Range("I132").Value = 2365.8667
$ws.Range("J132").Value = 3991.6365
$ws.Range("K132").Value = 7097.6001
$ws.Range("L132").Value = 11974.9095
$ws.Range("M132").Value = -4567.6001
$ws.Range("N132").Value = -17034.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2257.9167
$ws.Range("I7").Value = 2119.5
$ws.Range("J7").Value = 2950
$ws.Range("K7").Value = 2119.5
$ws.Range("L7").Value = 2950
$ws.Range("M7").Value = -2007.5
$ws.Range("N7").Value = -3174
$ws.Range("H40").Value = 2398.9033
$ws.Range("I40").Value = 2254.64
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2254.64
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -2118.64
$ws.Range("N40").Value = -3272
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H122").Value = 2713.95
$ws.Range("I122").Value = 2059.4
$ws.Range("K122").Value = 6178.200000000001
$ws.Range("M122").Value = -3728.200000000001
$ws.Range("H126").Value = 2257.9167
$ws.Range("I126").Value = 2119.5
$ws.Range("J126").Value = 2950
$ws.Range("K126").Value = 6358.5
$ws.Range("L126").Value = 8850
$ws.Range("M126").Value = -3888.5
$ws.Range("N126").Value = -13790
$ws.Range("H132").Value = 6989.5947
$ws.Range("I132").Value = 1400.2106
$ws.Range("J132").Value = 12889.5
$ws.Range("K132").Value = 4200.6318
$ws.Range("L132").Value = 38668.5
$ws.Range("M132").Value = -1670.6318
$ws.Range("N132").Value = -43728.5
$ws.Range("H136").Value = 5622.522
$ws.Range("I136").Value = 4193
$ws.Range("J136").Value = 7480.9
$ws.Range("K136").Value = 12579
$ws.Range("L136").Value = 22442.7
$ws.Range("M136").Value = -10029
$ws.Range("N136").Value = -27542.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 10000500
$ws.Range("I11").Value = 10000500
$ws.Range("K11").Value = 10000500
$ws.Range("M11").Value = -10000358
$ws.Range("H47").Value = 14339.667
$ws.Range("J47").Value = 14339.667
$ws.Range("L47").Value = 14339.667
$ws.Range("N47").Value = -15483.667
$ws.Range("H132").Value = 22149.66
$ws.Range("I132").Value = 35848.484
$ws.Range("J132").Value = 3232.238
$ws.Range("K132").Value = 107545.452
$ws.Range("L132").Value = 9696.714
$ws.Range("M132").Value = -105015.452
$ws.Range("N132").Value = -14756.714
$ws.Range("H136").Value = 25001902
$ws.Range("I136").Value = 43479660
$ws.Range("J136").Value = 2582.9412
$ws.Range("K136").Value = 130438980
$ws.Range("L136").Value = 7748.823600000001
$ws.Range("M136").Value = -130436430
$ws.Range("N136").Value = -12848.8236

